$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: new "Digi-Key Order 2" receipt line (barrel jacks order).
# D14 already carries the Hyperlink cell style; fill in the document name
# and wire up the external hyperlink to the new Digi-Key order 2 PDF.
$ws.Range("D14").Value = "Digi-Key Order 2.pdf"
$ws.Hyperlinks.Add($ws.Range("D14"), "Digi-Key%20Order%202.pdf", [Type]::Missing, [Type]::Missing, "Digi-Key Order 2.pdf")
# Re-apply the Hyperlink style so D14 keeps the same cell style as the
# other receipt-link cells instead of a freshly minted duplicate style.
$ws.Range("D14").Style = "Hyperlink"

# Cost of the new order.
$ws.Range("E14").Value = 21.51

# Selection ends up on D18 after the edits, matching the saved workbook view.
$ws.Range("D18").Select()
